# Add data for 2022-05-22 (through-date rolled from May 13 -> May 14 on the
# "current month" column B, plus the new day's carjacking counts bumped
# into the relevant neighborhood/month cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet tab and update the "current month" header label.
$ws.Name = "Through 2022-05-14"
$ws.Range("B1").Value = "May 2022 (through May 14)"

# Incremented existing counts.
$ws.Range("B2").Value = 4    # Englewood            / May 2022 (through May 14)
$ws.Range("L3").Value = 4    # Austin               / May 2020
$ws.Range("AA3").Value = 4   # Austin               / May 2017
$ws.Range("AF5").Value = 4   # Garfield Park        / May 2016
$ws.Range("V6").Value = 2    # Chicago Lawn         / May 2018
$ws.Range("B7").Value = 2    # North Lawndale       / May 2022 (through May 14)
$ws.Range("AF7").Value = 2   # North Lawndale       / May 2016
$ws.Range("G20").Value = 2   # Woodlawn             / May 2021
$ws.Range("G63").Value = 3   # Gage Park            / May 2021

# Newly populated (previously empty) cells.
$ws.Range("Q8").Value = 1    # South Shore          / May 2019
$ws.Range("AF8").Value = 1   # South Shore          / May 2016
$ws.Range("B17").Value = 1   # South Chicago        / May 2022 (through May 14)
$ws.Range("Q23").Value = 1   # Grand Crossing       / May 2019
$ws.Range("B24").Value = 1   # Grand Boulevard      / May 2022 (through May 14)
$ws.Range("Q51").Value = 1   # Ashburn              / May 2019
$ws.Range("Q53").Value = 1   # Boystown             / May 2019
$ws.Range("AF55").Value = 1  # Bucktown             / May 2016
$ws.Range("B57").Value = 1   # Chinatown            / May 2022 (through May 14)
$ws.Range("AF71").Value = 1  # Lincoln Square       / May 2016
$ws.Range("AA81").Value = 1  # Old Town             / May 2017
